# Adds a new "localdb" command-type column to the hidden '#system' sheet,
# which backs the named ranges used to validate Nexial command types/methods.
#
# Structural changes:
#   1. Insert a new column before column N ("macro"), shifting
#      macro..xml (N..AC) one column to the right (O..AD).
#   2. Populate the new column N with the "localdb" header and its six
#      method names (cloneTable, dropTables, exportCSV, importRecords,
#      purge, runSQLs).
#   3. Insert "localdb" into the alphabetical command-type list in column A
#      (row 14), shifting macro..xml down one row.
#   4. Re-point every defined name whose range moved, and add the new
#      "localdb" defined name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# --- 1. Insert the new column at N (14), pushing macro..xml to O..AD ---
$ws.Columns.Item(14).Insert()

# --- 2. Fill in the new "localdb" column ---
$ws.Cells.Item(1, 14).Value = "localdb"
$ws.Cells.Item(2, 14).Value = "cloneTable(var,source,target)"
$ws.Cells.Item(3, 14).Value = "dropTables(var,tables)"
$ws.Cells.Item(4, 14).Value = "exportCSV(sql,output)"
$ws.Cells.Item(5, 14).Value = "importRecords(var,sourceDb,sql,table)"
$ws.Cells.Item(6, 14).Value = "purge(var)"
$ws.Cells.Item(7, 14).Value = "runSQLs(var,sqls)"

# --- 3. Insert "localdb" alphabetically into the target list (column A) ---
# NOTE: Range.Insert() in this runtime always performs a full-row/column
# shift (like EntireRow/EntireColumn.Insert()), even for a single-cell
# range. Since only column A needs to move here, shift the A14:A29 values
# down to A15:A30 manually (bottom-up to avoid clobbering) instead.
for ($r = 29; $r -ge 14; $r--) {
    $ws.Cells.Item($r + 1, 1).Value = $ws.Cells.Item($r, 1).Value()
}
$ws.Cells.Item(14, 1).Value = "localdb"

# --- 4. Fix up defined names that referenced the shifted columns/rows ---
$names = $wb.Names

$names.Item("mail").RefersTo      = "='#system'!`$P`$2:`$P`$2"
$names.Item("number").RefersTo    = "='#system'!`$Q`$2:`$Q`$16"
$names.Item("pdf").RefersTo       = "='#system'!`$R`$2:`$R`$16"
$names.Item("rdbms").RefersTo     = "='#system'!`$S`$2:`$S`$7"
$names.Item("redis").RefersTo     = "='#system'!`$T`$2:`$T`$10"
$names.Item("sms").RefersTo       = "='#system'!`$U`$2:`$U`$2"
$names.Item("sound").RefersTo     = "='#system'!`$V`$2:`$V`$5"
$names.Item("ssh").RefersTo       = "='#system'!`$W`$2:`$W`$9"
$names.Item("step").RefersTo      = "='#system'!`$X`$2:`$X`$4"
$names.Item("target").RefersTo    = "='#system'!`$A`$2:`$A`$30"
$names.Item("web").RefersTo       = "='#system'!`$Y`$2:`$Y`$127"
$names.Item("webalert").RefersTo  = "='#system'!`$Z`$2:`$Z`$8"
$names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$names.Item("ws").RefersTo        = "='#system'!`$AB`$2:`$AB`$17"
$names.Item("ws.async").RefersTo  = "='#system'!`$AC`$2:`$AC`$8"
$names.Item("xml").RefersTo       = "='#system'!`$AD`$2:`$AD`$21"
$names.Item("macro").RefersTo     = "='#system'!`$O`$2:`$O`$4"

$names.Add("localdb", "='#system'!`$N`$2:`$N`$7")
